$d = $word.ActiveDocument

# 1. Remove the empty paragraph right after the "Auto Restart and Startup Settings" heading.
$d.Paragraphs.Item(2).Range.Delete()

# 2. Update the legacy VML picture's shape id and display size
#    (o:spid _x0000_i1027 -> _x0000_i1025; width:203.6pt;height:214.8pt -> width:204pt;height:214.5pt)
$picPara = $d.Paragraphs.Item(4)
$picRange = $picPara.Range
$shapeXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="5E1A0FA8" w14:textId="79BB001C" w:rsidR="001E46DF" w:rsidRDefault="00CA3F2C" w:rsidP="009404DC"><w:pPr><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana" w:cs="Verdana"/><w:lang w:val="en"/></w:rPr></w:pPr><w:r w:rsidRPr="00CA3F2C"><w:rPr><w:rFonts w:ascii="Verdana" w:hAnsi="Verdana" w:cs="Verdana"/><w:noProof/><w:lang w:val="en"/></w:rPr><w:pict w14:anchorId="60F96502"><v:shapetype id="_x0000_t75" coordsize="21600,21600" o:spt="75" o:preferrelative="t" path="m@4@5l@4@11@9@11@9@5xe" filled="f" stroked="f"><v:stroke joinstyle="miter"/><v:formulas><v:f eqn="if lineDrawn pixelLineWidth 0"/><v:f eqn="sum @0 1 0"/><v:f eqn="sum 0 0 @1"/><v:f eqn="prod @2 1 2"/><v:f eqn="prod @3 21600 pixelWidth"/><v:f eqn="prod @3 21600 pixelHeight"/><v:f eqn="sum @0 0 1"/><v:f eqn="prod @6 1 2"/><v:f eqn="prod @7 21600 pixelWidth"/><v:f eqn="sum @8 21600 0"/><v:f eqn="prod @7 21600 pixelHeight"/><v:f eqn="sum @10 21600 0"/></v:formulas><v:path o:extrusionok="f" gradientshapeok="t" o:connecttype="rect"/><o:lock v:ext="edit" aspectratio="t"/></v:shapetype><v:shape id="Picture 1" o:spid="_x0000_i1025" type="#_x0000_t75" style="width:204pt;height:214.5pt;visibility:visible;mso-wrap-style:square"><v:imagedata r:id="rId5" o:title=""/></v:shape></w:pict></w:r></w:p>
'@
$picRange.InsertXML($shapeXml)
